$d = $word.ActiveDocument
$wordNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Locate the existing (empty) "_GoBack" bookmark paragraph - it is the
# paragraph immediately preceding the "Description (4000 chars)" Heading1
# paragraph.
$bookmarkPara = $null
foreach ($p in $d.Paragraphs) {
    $next = $p.Next()
    if ($next -ne $null -and $next.Range.Text.TrimEnd() -eq "Description (4000 chars)") {
        $bookmarkPara = $p
    }
}

$bmStart = $bookmarkPara.Range.Start

# 1) Insert a brand new, empty Heading1 paragraph right before the bookmark
#    paragraph.
$r1 = $d.Range($bmStart, $bmStart)
$r1.InsertXML('<w:p ' + $wordNs + '><w:pPr><w:pStyle w:val="Heading1"/></w:pPr></w:p>')

# Recompute the bookmark paragraph start (it moved forward).
$bookmarkPara = $null
foreach ($p in $d.Paragraphs) {
    $next = $p.Next()
    if ($next -ne $null -and $next.Range.Text.TrimEnd() -eq "Description (4000 chars)") {
        $bookmarkPara = $p
    }
}
$bmStart = $bookmarkPara.Range.Start

# 2) Give the bookmark paragraph the Heading1 style.
$bookmarkPara.Style = "Heading1"

# 3) Insert the "Short Description (80 cha" runs immediately before the
#    (still empty) bookmark, i.e. at the very start of the paragraph.
$r2 = $d.Range($bmStart, $bmStart)
$r2.InsertXML('<w:p ' + $wordNs + '>' +
  '<w:r><w:t xml:space="preserve">Short </w:t></w:r>' +
  '<w:r><w:t>Description (</w:t></w:r>' +
  '<w:r><w:t>8</w:t></w:r>' +
  '<w:r><w:t>0 cha</w:t></w:r>' +
  '</w:p>')

# 4) Insert "rs)" right after the bookmark (still within the same
#    paragraph). The bookmark now sits right after the text just inserted.
$afterBookmark = $bmStart + ("Short Description (80 cha").Length
$r3 = $d.Range($afterBookmark, $afterBookmark)
$r3.InsertAfter("rs)")

# 5) Insert the new "Kidoju provides..." paragraph and the following blank
#    paragraph right before the "Description (4000 chars)" heading. Both
#    paragraphs must be supplied in a single InsertXML call: when the
#    target range sits at the start of a non-empty paragraph, only the
#    *last* <w:p> of the payload is merged into that destination
#    paragraph, while any earlier <w:p> elements become genuine new
#    paragraphs of their own.
$descHeading = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "Description (4000 chars)") {
        $descHeading = $p
    }
}
$dhStart = $descHeading.Range.Start
$r4 = $d.Range($dhStart, $dhStart)
$r4.InsertXML('<w:p ' + $wordNs + '><w:r><w:t>Kidoju provides many auto-corrected exercises to help children do their best in school tests.</w:t></w:r></w:p>' +
  '<w:p ' + $wordNs + '></w:p>')
